$wb = $excel.ActiveWorkbook

# Row 17 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 10441991
$ws.Range("J17").Value = 20883332
$ws.Range("L17").Value = 62649996
$ws.Range("N17").Value = -62650332

# Row 33 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2100.7144
$ws.Range("I33").Value = 77.86667
$ws.Range("K33").Value = 77.86667
$ws.Range("M33").Value = 151.13333

# Row 41 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 5801.1816
$ws.Range("I41").Value = 2612.5
$ws.Range("J41").Value = 7623.2856
$ws.Range("K41").Value = 2612.5
$ws.Range("L41").Value = 7623.2856
$ws.Range("M41").Value = -2172.5
$ws.Range("N41").Value = -8503.285599999999

# Row 137 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3330.9058
$ws.Range("I137").Value = 2573.4644
$ws.Range("J137").Value = 4179.24
$ws.Range("K137").Value = 7720.3932
$ws.Range("L137").Value = 12537.72
$ws.Range("M137").Value = -5170.3932
$ws.Range("N137").Value = -17637.72

# Row 138 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5062.375
$ws.Range("J138").Value = 6869.6
$ws.Range("L138").Value = 20608.8
$ws.Range("N138").Value = -30888.8

# Row 74 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7939541.5
$ws.Range("I74").Value = 9806287
$ws.Range("J74").Value = 5875.875
$ws.Range("K74").Value = 9806287
$ws.Range("L74").Value = 5875.875
$ws.Range("M74").Value = -9805413
$ws.Range("N74").Value = -7623.875

# Row 77 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7939541.5
$ws.Range("I77").Value = 9806287
$ws.Range("J77").Value = 5875.875
$ws.Range("K77").Value = 49031435
$ws.Range("L77").Value = 29379.375
$ws.Range("M77").Value = -49027067
$ws.Range("N77").Value = -38115.375

# Row 102 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3435.1667
$ws.Range("J102").Value = 3393
$ws.Range("L102").Value = 3393
$ws.Range("N102").Value = -6637

# Row 20 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2217
$ws.Range("I20").Value = 1564.4
$ws.Range("K20").Value = 1564.4
$ws.Range("M20").Value = -1317.4

# Row 99 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3305.4285
$ws.Range("I99").Value = 2916.5264
$ws.Range("K99").Value = 2916.5264
$ws.Range("M99").Value = -1418.5264

# Row 134 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3820.2
$ws.Range("I134").Value = 2189.3845
$ws.Range("J134").Value = 6848.857
$ws.Range("K134").Value = 6568.1535
$ws.Range("L134").Value = 20546.571
$ws.Range("M134").Value = -4033.1535
$ws.Range("N134").Value = -25616.571

# Row 16 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2733.3333
$ws.Range("I16").Value = 1832.25
$ws.Range("J16").Value = 6337.6665
$ws.Range("K16").Value = 1832.25
$ws.Range("L16").Value = 6337.6665
$ws.Range("M16").Value = -1545.25
$ws.Range("N16").Value = -6911.6665

# Row 22 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1313.7916
$ws.Range("J22").Value = 4090
$ws.Range("L22").Value = 4090
$ws.Range("N22").Value = -4790

# Row 31 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19175.94
$ws.Range("I31").Value = 1742
$ws.Range("K31").Value = 1742
$ws.Range("M31").Value = -1447

# Row 34 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 19175.94
$ws.Range("I34").Value = 1742
$ws.Range("K34").Value = 1742
$ws.Range("M34").Value = -1540

# Row 94 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4804.1665
$ws.Range("I94").Value = 1937.3334
$ws.Range("K94").Value = 1937.3334
$ws.Range("M94").Value = -1486.3334

# Row 99 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3087.5
$ws.Range("I99").Value = 3039.1
$ws.Range("K99").Value = 3039.1
$ws.Range("M99").Value = -1541.1

# Row 113 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2733.3333
$ws.Range("I113").Value = 1832.25
$ws.Range("J113").Value = 6337.6665
$ws.Range("K113").Value = 1832.25
$ws.Range("L113").Value = 6337.6665
$ws.Range("M113").Value = 337.75
$ws.Range("N113").Value = -10677.6665

# Row 126 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3087.5
$ws.Range("I126").Value = 3039.1
$ws.Range("K126").Value = 9117.299999999999
$ws.Range("M126").Value = -6647.299999999999

# Row 132 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3955.5186
$ws.Range("I132").Value = 2764.7058
$ws.Range("J132").Value = 5979.9
$ws.Range("K132").Value = 8294.117400000001
$ws.Range("L132").Value = 17939.7
$ws.Range("M132").Value = -5764.117400000001
$ws.Range("N132").Value = -22999.7

# Row 134 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2957.3684
$ws.Range("I134").Value = 2065.3333
$ws.Range("K134").Value = 6195.999899999999
$ws.Range("M134").Value = -3660.999899999999

# Row 137 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1114385.9
$ws.Range("I137").Value = 1667240
$ws.Range("J137").Value = 8677.666999999999
$ws.Range("K137").Value = 5001720
$ws.Range("L137").Value = 26033.001
$ws.Range("M137").Value = -4996620
$ws.Range("N137").Value = -36233.001

# Row 14 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1821.4286
$ws.Range("I14").Value = 333.33334
$ws.Range("J14").Value = 10750
$ws.Range("K14").Value = 333.33334
$ws.Range("L14").Value = 10750
$ws.Range("M14").Value = -165.33334
$ws.Range("N14").Value = -11086

# Row 132 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6837.95
$ws.Range("I132").Value = 4978.5557
$ws.Range("J132").Value = 8359.272000000001
$ws.Range("K132").Value = 14935.6671
$ws.Range("L132").Value = 25077.816
$ws.Range("M132").Value = -12405.6671
$ws.Range("N132").Value = -30137.816

# Row 46 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3442.9092
$ws.Range("J46").Value = 5083.6665
$ws.Range("L46").Value = 5083.6665
$ws.Range("N46").Value = -5459.6665

# Row 93 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 11373.469
$ws.Range("I93").Value = 11294.833
$ws.Range("K93").Value = 11294.833
$ws.Range("M93").Value = -10046.833

# Row 122 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 23176.928
$ws.Range("I122").Value = 26747.7
$ws.Range("J122").Value = 14250
$ws.Range("K122").Value = 80243.10000000001
$ws.Range("L122").Value = 42750
$ws.Range("M122").Value = -77793.10000000001
$ws.Range("N122").Value = -47650

# Row 132 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4426.8237
$ws.Range("I132").Value = 3318.9524
$ws.Range("J132").Value = 6216.4614
$ws.Range("K132").Value = 9956.8572
$ws.Range("L132").Value = 18649.3842
$ws.Range("M132").Value = -7426.8572
$ws.Range("N132").Value = -23709.3842

# Row 136 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9008.421
$ws.Range("I136").Value = 4268.75
$ws.Range("J136").Value = 12455.454
$ws.Range("K136").Value = 12806.25
$ws.Range("L136").Value = 37366.362
$ws.Range("M136").Value = -10256.25
$ws.Range("N136").Value = -42466.362

# Row 15 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 133403.1
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 133403.1
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 133403.1
$ws.Range("M15").Value = $null
$ws.Range("N15").Value = -133979.1

# Row 38 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 3062
$ws.Range("J38").Value = 3062
$ws.Range("L38").Value = 3062
$ws.Range("M38").Value = -4008

# Row 132 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6062.162
$ws.Range("I132").Value = 2816
$ws.Range("K132").Value = 8448
$ws.Range("M132").Value = -5918

# Row 136 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4012.2927
$ws.Range("I136").Value = 3584
$ws.Range("J136").Value = 6510.6665
$ws.Range("K136").Value = 10752
$ws.Range("L136").Value = 19531.9995
$ws.Range("M136").Value = -8202
